$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the tiny floating-point precision of the existing A14 timestamp
$ws.Range("A14").Value = 45814.39350229167

# Append the new data row (row 15)
$ws.Range("A15").Value = 45815.3911123033
$ws.Range("A15").NumberFormat = $ws.Range("A14").NumberFormat

$ws.Range("B15").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C15").Value = "1Kg"
$ws.Range("D15").Value = "15,41€"
